# Update gh-pages output data (regenerated scrape values) for 杭州-漫展信息.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 6242
$ws1.Range("F6").Value  = 1809
$ws1.Range("F7").Value  = 6378
$ws1.Range("F11").Value = 15
$ws1.Range("F17").Value = 7993
$ws1.Range("F21").Value = 110
$ws1.Range("F22").Value = 1750
$ws1.Range("F29").Value = 1799
$ws1.Range("F31").Value = 384
$ws1.Range("F34").Value = 11
$ws1.Range("F35").Value = 86
$ws1.Range("F37").Value = 3925

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9544
$ws3.Range("F3").Value = 2284
$ws3.Range("F4").Value = 691

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 9544
$ws4.Range("F3").Value  = 2284
$ws4.Range("F4").Value  = 691
$ws4.Range("F10").Value = 6242
$ws4.Range("F12").Value = 1809
$ws4.Range("F13").Value = 6378
$ws4.Range("F24").Value = 7993
$ws4.Range("F27").Value = 110
$ws4.Range("F28").Value = 1750
$ws4.Range("F32").Value = 1799
$ws4.Range("F39").Value = 11
$ws4.Range("F44").Value = 3925

$wb.Save()
